# Trade #45 closed at 2026-02-17 21:07:22 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#  - Summary sheet totals refreshed
#  - Strategy Status row for MarketMaking refreshed
#  - "All Trades" + "MarketMaking" sheets: trade #73 marked CLOSED (early_exit)
#    and a brand-new open trade #106 appended to both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.5    # Current Capital
$summary.Range("B4").Value = 0.3       # Total P&L $
$summary.Range("B5").Value = 0.08      # Total P&L %
$summary.Range("B6").Value = 73        # Total Trades
$summary.Range("B8").Value = 30        # Losing Trades
$summary.Range("B9").Value = 45.21     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.5      # Capital
$status.Range("D5").Value = 40         # Trades
$status.Range("E5").Value = 0.19       # P&L $
$status.Range("F5").Value = 0.5        # P&L %
$status.Range("G5").Value = 47.5       # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet - close out trade #73 (row 74)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G74").Value = 0.16
$allTrades.Range("H74").Value = "CLOSED"
$allTrades.Range("I74").Value = -23.8095
$allTrades.Range("J74").Value = -0.05
$allTrades.Range("K74").Value = 100.5
$allTrades.Range("L74").Value = "early_exit"
$allTrades.Range("M74").Value = 0.13

# New trade #106 appended as row 107
$allTrades.Cells.Item(107, 1).Value = 106
$allTrades.Cells.Item(107, 2).NumberFormat = "@"
$allTrades.Cells.Item(107, 2).Value = "2026-02-17"
$allTrades.Cells.Item(107, 2).Style = "Normal"
$allTrades.Cells.Item(107, 3).NumberFormat = "@"
$allTrades.Cells.Item(107, 3).Value = "21:07:15"
$allTrades.Cells.Item(107, 3).Style = "Normal"
$allTrades.Cells.Item(107, 4).Value = "MarketMaking"
$allTrades.Cells.Item(107, 5).Value = "DOWN"
$allTrades.Cells.Item(107, 6).Value = 0.21
$allTrades.Cells.Item(107, 8).Value = "OPEN"
$allTrades.Cells.Item(107, 9).Value = 0
$allTrades.Cells.Item(107, 10).Value = 0
$allTrades.Cells.Item(107, 11).Value = 100.5519219857093
$allTrades.Cells.Item(107, 13).Value = 0
$allTrades.Cells.Item(107, 14).Value = 0
$allTrades.Cells.Item(107, 15).Value = 0
$allTrades.Cells.Item(107, 16).Value = 0.6
$allTrades.Cells.Item(107, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet - close out trade #73 (row 41)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G41").Value = 0.16
$mm.Range("H41").Value = "CLOSED"
$mm.Range("I41").Value = -23.8095
$mm.Range("J41").Value = -0.05
$mm.Range("K41").Value = 100.5
$mm.Range("P41").Value = "early_exit"
$mm.Range("Q41").Value = 0.13

# New trade #106 appended as row 74
$mm.Cells.Item(74, 1).Value = 106
$mm.Cells.Item(74, 2).NumberFormat = "@"
$mm.Cells.Item(74, 2).Value = "2026-02-17"
$mm.Cells.Item(74, 2).Style = "Normal"
$mm.Cells.Item(74, 3).NumberFormat = "@"
$mm.Cells.Item(74, 3).Value = "21:07:15"
$mm.Cells.Item(74, 3).Style = "Normal"
$mm.Cells.Item(74, 4).Value = "MarketMaking"
$mm.Cells.Item(74, 5).Value = "DOWN"
$mm.Cells.Item(74, 6).Value = 0.21
$mm.Cells.Item(74, 8).Value = "OPEN"
$mm.Cells.Item(74, 9).Value = 0
$mm.Cells.Item(74, 10).Value = 0
$mm.Cells.Item(74, 11).Value = 100.5519219857093
$mm.Cells.Item(74, 12).Value = 0
$mm.Cells.Item(74, 13).Value = 0
$mm.Cells.Item(74, 14).Value = 0.6
$mm.Cells.Item(74, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(74, 17).Value = 0
